# analysis_draft2017b: updating analysis to use marc_s2 simulation data
# Update the dust/sea-salt aerosol statistics table (B2:C8) with new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sources, Tg/yr
$ws.Range("B2").Value = "+3683.19 ± 25.59"
$ws.Range("C2").Value = "+5484.88 ± 11.18"

#     Emission
$ws.Range("B3").Value = "+3683.19 ± 25.59"
$ws.Range("C3").Value = "+5484.88 ± 11.18"

# Sinks, Tg/yr
$ws.Range("B4").Value = "-3705.79 ± 25.88"
$ws.Range("C4").Value = "-5533.62 ± 11.24"

#     Impaction scavenging
$ws.Range("B5").Value = "-1819.27 ± 11.26"
$ws.Range("C5").Value = "-2324.42 ± 4.29"

#     Dry deposition
$ws.Range("B6").Value = "-1886.52 ± 16.05"
$ws.Range("C6").Value = "-3209.19 ± 7.27"

# Burden, Tg
$ws.Range("B7").Value = "+40.91 ± 0.30"
$ws.Range("C7").Value = "+9.60 ± 0.02"

# Lifetime, days (B8 changes, C8 is unchanged: "+0.63 ± 0.00")
$ws.Range("B8").Value = "+4.03 ± 0.04"
